$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 3.052282377742472
$ws.Range("C2").Value = 0.1354979047511904
$ws.Range("D2").Value = 0.4263785751512756
$ws.Range("E2").Value = 0.1119939907509493
$ws.Range("G2").Value = 3.108132831540189
$ws.Range("H2").Value = 2.388476113548734
$ws.Range("J2").Value = 0.03700484573729312
$ws.Range("L2").Value = 0.4786374334240406
$ws.Range("M2").Value = 0.6742196960079312
$ws.Range("N2").Value = 2.74755183029724
$ws.Range("B3").Value = 2.965402681747264
$ws.Range("C3").Value = 0.1187811219180617
$ws.Range("D3").Value = 0.4265951139999942
$ws.Range("E3").Value = 0.1124319677617645
$ws.Range("G3").Value = 3.093615699130311
$ws.Range("H3").Value = 2.389290011603123
$ws.Range("J3").Value = 0.03554302740950277
$ws.Range("L3").Value = 0.4762238575677173
$ws.Range("M3").Value = 0.6608947510589971
$ws.Range("N3").Value = 2.771258202610028
$ws.Range("B4").Value = 2.913729657881333
$ws.Range("C4").Value = 0.1085397881158769
$ws.Range("D4").Value = 0.4268870576316317
$ws.Range("E4").Value = 0.112718012275093
$ws.Range("G4").Value = 3.086350049555591
$ws.Range("H4").Value = 2.390846973965779
$ws.Range("J4").Value = 0.03463242925914756
$ws.Range("L4").Value = 0.4749467923772883
$ws.Range("M4").Value = 0.6530646173872725
$ws.Range("N4").Value = 2.78661955006266
$ws.Range("B5").Value = 2.893092506911159
$ws.Range("C5").Value = 0.1043719279981019
$ws.Range("D5").Value = 0.427046035681883
$ws.Range("E5").Value = 0.1128388930679805
$ws.Range("G5").Value = 3.083802651474713
$ws.Range("H5").Value = 2.391747028020148
$ws.Range("J5").Value = 0.03425805048743058
$ws.Range("L5").Value = 0.4744779276586826
$ws.Range("M5").Value = 0.6499621260447412
$ws.Range("N5").Value = 2.79308176011417
$ws.Range("B6").Value = 2.889691084813137
$ws.Range("C6").Value = 0.1036801891242476
$ws.Range("D6").Value = 0.4270748511931046
$ws.Range("E6").Value = 0.1128592261363619
$ws.Range("G6").Value = 3.083404598971015
$ws.Range("H6").Value = 2.391912513948085
$ws.Range("J6").Value = 0.03419568498811643
$ws.Range("L6").Value = 0.4744031874609504
$ws.Range("M6").Value = 0.6494522954870519
$ws.Range("N6").Value = 2.794167019334218
$ws.Range("B7").Value = 2.913449637741564
$ws.Range("C7").Value = 0.1084835565863784
$ws.Range("D7").Value = 0.4268890396313694
$ws.Range("E7").Value = 0.1127196250309139
$ws.Range("G7").Value = 3.086314021857419
$ws.Range("H7").Value = 2.390858037463488
$ws.Range("J7").Value = 0.03462739366023726
$ws.Range("L7").Value = 0.474940260338677
$ws.Range("M7").Value = 0.6530224183965956
$ws.Range("N7").Value = 2.786705882750184
$ws.Range("B8").Value = 3.021979138070549
$ws.Range("C8").Value = 0.1297291060336079
$ws.Range("D8").Value = 0.4264202566036914
$ws.Range("E8").Value = 0.112141456467227
$ws.Range("G8").Value = 3.102784695945445
$ws.Range("H8").Value = 2.388537149288652
$ws.Range("J8").Value = 0.0365035011367425
$ws.Range("L8").Value = 0.47776272575517
$ws.Range("M8").Value = 0.6695522951634274
$ws.Range("N8").Value = 2.755558465034916
$ws.Range("B9").Value = 3.248096114968291
$ws.Range("C9").Value = 0.1715828994540516
$ws.Range("D9").Value = 0.4267615777836369
$ws.Range("E9").Value = 0.1111431341681439
$ws.Range("G9").Value = 3.148209005755405
$ws.Range("H9").Value = 2.392389392350225
$ws.Range("J9").Value = 0.04008036029789608
$ws.Range("L9").Value = 0.4849227491686179
$ws.Range("M9").Value = 0.7047594950828397
$ws.Range("N9").Value = 2.700879314503133
$ws.Range("B10").Value = 3.422389875593694
$ws.Range("C10").Value = 0.2024666677810671
$ws.Range("D10").Value = 0.4277802629171958
$ws.Range("E10").Value = 0.1104916761912067
$ws.Range("G10").Value = 3.189661388756406
$ws.Range("H10").Value = 2.400367256814491
$ws.Range("J10").Value = 0.04264791690571457
$ws.Range("L10").Value = 0.4911747046135559
$ws.Range("M10").Value = 0.7323375021413838
$ws.Range("N10").Value = 2.664620316393858
$ws.Range("B11").Value = 3.503469181791161
$ws.Range("C11").Value = 0.2165494670076953
$ws.Range("D11").Value = 0.4284103762338418
$ws.Range("E11").Value = 0.1102129975332158
$ws.Range("G11").Value = 3.210290985972193
$ws.Range("H11").Value = 2.405120079281346
$ws.Range("J11").Value = 0.04380328639687647
$ws.Range("L11").Value = 0.4942344556279608
$ws.Range("M11").Value = 0.745257331686382
$ws.Range("N11").Value = 2.648977579256396
$ws.Range("B12").Value = 3.534430353566165
$ws.Range("C12").Value = 0.2218873459271435
$ws.Range("D12").Value = 0.4286729434231091
$ws.Range("E12").Value = 0.1101100013685814
$ws.Range("G12").Value = 3.218359087717545
$ws.Range("H12").Value = 2.407081848142525
$ws.Range("J12").Value = 0.04423900948466297
$ws.Range("L12").Value = 0.4954241241755568
$ws.Range("M12").Value = 0.750203696967283
$ws.Range("N12").Value = 2.643176820268096
$ws.Range("B13").Value = 3.527750819820653
$ws.Range("C13").Value = 0.2207375114201398
$ws.Range("D13").Value = 0.4286153296225734
$ws.Range("E13").Value = 0.110132070887687
$ws.Range("G13").Value = 3.216610067137736
$ws.Range("H13").Value = 2.406652135634261
$ws.Range("J13").Value = 0.04414524794096764
$ws.Range("L13").Value = 0.4951665288617875
$ws.Range("M13").Value = 0.7491360101940643
$ws.Range("N13").Value = 2.644420650205205
$ws.Range("B14").Value = 3.506011198482327
$ws.Range("C14").Value = 0.2169885156141902
$ws.Range("D14").Value = 0.4284314977458052
$ws.Range("E14").Value = 0.1102044732489227
$ws.Range("G14").Value = 3.210949613519915
$ws.Range("H14").Value = 2.405278226584869
$ws.Range("J14").Value = 0.04383916933254639
$ws.Range("L14").Value = 0.4943317089982031
$ws.Range("M14").Value = 0.7456631914109053
$ws.Range("N14").Value = 2.648497884145563
$ws.Range("B15").Value = 3.492728692539174
$ws.Range("C15").Value = 0.2146928081191959
$ws.Range("D15").Value = 0.4283220148786029
$ws.Range("E15").Value = 0.1102491514744424
$ws.Range("G15").Value = 3.207515812957809
$ws.Range("H15").Value = 2.404457773955443
$ws.Range("J15").Value = 0.04365145495191314
$ws.Range("L15").Value = 0.4938243955015054
$ws.Range("M15").Value = 0.7435430116861781
$ws.Range("N15").Value = 2.651011310790956
$ws.Range("B16").Value = 3.417127256410083
$ws.Range("C16").Value = 0.2015470167361286
$ws.Range("D16").Value = 0.4277424367723626
$ws.Range("E16").Value = 0.1105102434391907
$ws.Range("G16").Value = 3.188348970218243
$ws.Range("H16").Value = 2.400079293251082
$ws.Range("J16").Value = 0.04257215889546373
$ws.Range("L16").Value = 0.4909790822438822
$ws.Range("M16").Value = 0.731500700902032
$ws.Range("N16").Value = 2.665659772411821
$ws.Range("B17").Value = 3.371207582880515
$ws.Range("C17").Value = 0.1934912469836547
$ws.Range("D17").Value = 0.4274295709301015
$ws.Range("E17").Value = 0.1106749356669012
$ws.Range("G17").Value = 3.177045618071645
$ws.Range("H17").Value = 2.397681311482785
$ws.Range("J17").Value = 0.04190683004630102
$ws.Range("L17").Value = 0.4892888135009201
$ws.Range("M17").Value = 0.7242090813041102
$ws.Range("N17").Value = 2.674864498579929
$ws.Range("B18").Value = 3.344964460913843
$ws.Range("C18").Value = 0.1888609337467528
$ws.Range("D18").Value = 0.4272653109993172
$ws.Range("E18").Value = 0.1107713262016508
$ws.Range("G18").Value = 3.170711010405341
$ws.Range("H18").Value = 2.396407794800268
$ws.Range("J18").Value = 0.0415229606587495
$ws.Range("L18").Value = 0.4883369193264286
$ws.Range("M18").Value = 0.7200503762031119
$ws.Range("N18").Value = 2.68023895552745
$ws.Range("B19").Value = 3.33610794399857
$ws.Range("C19").Value = 0.187293724824741
$ws.Range("D19").Value = 0.4272123910568979
$ws.Range("E19").Value = 0.1108042484038223
$ws.Range("G19").Value = 3.168594829442526
$ws.Range("H19").Value = 2.395994752725187
$ws.Range("J19").Value = 0.04139278396132084
$ws.Range("L19").Value = 0.4880181118628997
$ws.Range("M19").Value = 0.7186483611579035
$ws.Range("N19").Value = 2.682072408561403
$ws.Range("B20").Value = 3.376078351519311
$ws.Range("C20").Value = 0.194348469756676
$ws.Range("D20").Value = 0.427461252064262
$ws.Range("E20").Value = 0.1106572317464254
$ws.Range("G20").Value = 3.178231609301434
$ws.Range("H20").Value = 2.397925634345654
$ws.Range("J20").Value = 0.04197777839837258
$ws.Range("L20").Value = 0.4894666442262263
$ws.Range("M20").Value = 0.7249816388722508
$ws.Range("N20").Value = 2.673876343002959
$ws.Range("B21").Value = 3.51238964082836
$ws.Range("C21").Value = 0.2180895485207657
$ws.Range("D21").Value = 0.4284848435560349
$ws.Range("E21").Value = 0.1101831382276908
$ws.Range("G21").Value = 3.212605265176819
$ws.Range("H21").Value = 2.405677377708685
$ws.Range("J21").Value = 0.04392912042494501
$ws.Range("L21").Value = 0.4945760743761696
$ws.Range("M21").Value = 0.7466817789650122
$ws.Range("N21").Value = 2.64729696605059
$ws.Range("B22").Value = 3.602982170275368
$ws.Range("C22").Value = 0.2336351185796843
$ws.Range("D22").Value = 0.4292934529657373
$ws.Range("E22").Value = 0.1098880536475644
$ws.Range("G22").Value = 3.236563880519128
$ws.Range("H22").Value = 2.411687915078346
$ws.Range("J22").Value = 0.04519401202984596
$ws.Range("L22").Value = 0.498096112884042
$ws.Range("M22").Value = 0.7611783134974104
$ws.Range("N22").Value = 2.630641776611917
$ws.Range("B23").Value = 3.554493337238966
$ws.Range("C23").Value = 0.2253353997895999
$ws.Range("D23").Value = 0.4288491112341859
$ws.Range("E23").Value = 0.110044197632206
$ws.Range("G23").Value = 3.223639661842384
$ws.Range("H23").Value = 2.408393437983023
$ws.Range("J23").Value = 0.04451986093283722
$ws.Range("L23").Value = 0.4962008682239372
$ws.Range("M23").Value = 0.7534124660916603
$ws.Range("N23").Value = 2.639465339757898
$ws.Range("B24").Value = 3.373875788919065
$ws.Range("C24").Value = 0.1939609160696421
$ws.Range("D24").Value = 0.4274468803983211
$ws.Range("E24").Value = 0.1106652303720304
$ws.Range("G24").Value = 3.177694912408924
$ws.Range("H24").Value = 2.397814848574825
$ws.Range("J24").Value = 0.04194570689613997
$ws.Range("L24").Value = 0.48938618508204
$ws.Range("M24").Value = 0.7246322617374474
$ws.Range("N24").Value = 2.674322831056521
$ws.Range("B25").Value = 3.185496178490212
$ws.Range("C25").Value = 0.1602380961002723
$ws.Range("D25").Value = 0.4265343499005638
$ws.Range("E25").Value = 0.1113987639712448
$ws.Range("G25").Value = 3.134507414773935
$ws.Range("H25").Value = 2.390445097700621
$ws.Range("J25").Value = 0.03912346538349709
$ws.Range("L25").Value = 0.4828117189441343
$ws.Range("M25").Value = 0.694935173350288
$ws.Range("N25").Value = 2.714985042102938
